$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly crypto price/volume refresh (cell-by-cell update per upstream diff).
# Numeric-looking Price (column D) values must stay as text, matching the
# original inline-string cells, so they are entered with a leading apostrophe
# (the same trick a user would use in the Excel UI) to avoid Excel silently
# converting them to numbers.

# Row 2
$ws.Range("D2").Value = "33.666.04"
$ws.Range("E2").Value = "  -0.59%  "

# Row 3
$ws.Range("D3").Value = "1.763.77"
$ws.Range("E3").Value = "  -0.88%  "

# Row 4
$ws.Range("E4").Value = "  +0.35%  "

# Row 5
$ws.Range("D5").Value = "'224.22"
$ws.Range("E5").Value = "  +1.42%  "

# Row 6
$ws.Range("E6").Value = "  -1.42%  "

# Row 7
$ws.Range("E7").Value = "  +0.33%  "

# Row 8
$ws.Range("D8").Value = "'31.87"
$ws.Range("E8").Value = "  +2.61%  "

# Row 9
$ws.Range("E9").Value = "  +0.69%  "

# Row 10
$ws.Range("D10").Value = "'0.0685"
$ws.Range("E10").Value = "  -3.35%  "

# Row 12
$ws.Range("D12").Value = "2.019.40"
$ws.Range("E12").Value = "  -0.70%  "

# Row 13
$ws.Range("D13").Value = "'11.22"
$ws.Range("E13").Value = "  +7.12%  "

# Row 14
$ws.Range("D14").Value = "1.754.20"
$ws.Range("E14").Value = "  -1.42%  "

# Row 15
$ws.Range("D15").Value = "33.695.05"
$ws.Range("E15").Value = "  -0.49%  "

# Row 16
$ws.Range("E16").Value = "  -2.45%  "

# Row 17
$ws.Range("D17").Value = "'4.12"
$ws.Range("E17").Value = "  -2.11%  "

# Row 18
$ws.Range("D18").Value = "'66.51"
$ws.Range("E18").Value = "  -2.14%  "

# Row 19
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.0₃0770"
$ws.Range("E19").Value = "  -0.71%  "

# Row 20
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "'237.04"
$ws.Range("E20").Value = "  -3.14%  "

# Row 21
$ws.Range("E21").Value = "  +0.22%  "

# Row 22
$ws.Range("D22").Value = "'10.57"
$ws.Range("E22").Value = "  -0.62%  "

# Row 23
$ws.Range("D23").Value = "'4.03"
$ws.Range("E23").Value = "  -1.06%  "

# Row 24
$ws.Range("E24").Value = "  -1.81%  "

# Row 25
$ws.Range("D25").Value = "'159.54"
$ws.Range("E25").Value = "  +1.35%  "

# Row 26
$ws.Range("D26").Value = "'16.09"
$ws.Range("E26").Value = "  -1.77%  "

# Row 27
$ws.Range("D27").Value = "'7.02"
$ws.Range("E27").Value = "  +0.63%  "

# Row 28
$ws.Range("E28").Value = "  -0.23%  "

# Row 29
$ws.Range("E29").Value = "  +0.49%  "

# Row 30
$ws.Range("E30").Value = "  +2.21%  "

# Row 31
$ws.Range("D31").Value = "'0.0509"
$ws.Range("E31").Value = "  -2.08%  "

# Row 32
$ws.Range("E32").Value = "  -3.04%  "

# Row 33
$ws.Range("E33").Value = "  -0.01%  "

# Row 34
$ws.Range("E34").Value = "  -1.84%  "

# Row 35
$ws.Range("D35").Value = "1.377.91"
$ws.Range("E35").Value = "  -1.30%  "

# Row 36
$ws.Range("E36").Value = "  +2.29%  "

# Row 37
$ws.Range("E37").Value = "  -1.50%  "

# Row 38
$ws.Range("D38").Value = "'0.0183"
$ws.Range("E38").Value = "  -1.13%  "

# Row 39
$ws.Range("D39").Value = "'2.21"
$ws.Range("E39").Value = "  +4.77%  "

# Row 40
$ws.Range("E40").Value = "  +0.70%  "

# Row 41
$ws.Range("E41").Value = "  -2.79%  "

# Row 42
$ws.Range("D42").Value = "'77.51"
$ws.Range("E42").Value = "  -1.92%  "

# Row 43
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").Value = "'2.65"
$ws.Range("E43").Value = "  -1.79%  "

# Row 44
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").Value = "'13.47"
$ws.Range("E44").Value = "  +15.04%  "

# Row 45
$ws.Range("D45").Value = "0.0₆0140"
$ws.Range("E45").Value = "  +16.64%  "

# Row 46
$ws.Range("E46").Value = "  +4.61%  "

# Row 47
$ws.Range("D47").Value = "'0.0500"
$ws.Range("E47").Value = "  +2.22%  "

# Row 48
$ws.Range("D48").Value = "'107.60"
$ws.Range("E48").Value = "  +2.61%  "

# Row 49
$ws.Range("E49").Value = "  -1.89%  "

# Row 50
$ws.Range("D50").Value = "1.919.92"
$ws.Range("E50").Value = "  +0.00%  "

# Row 51
$ws.Range("E51").Value = "  +0.61%  "
